$d = $word.ActiveDocument

$replacements = @(
    @{old="2025-11-21 Friday"; new="2025-11-22 Saturday"},
    @{old="335÷8=41, 7"; new="431÷4=107, 3"},
    @{old="307÷7=43, 6"; new="572÷4=143, 0"},
    @{old="166÷4=41, 2"; new="437÷8=54, 5"},
    @{old="678÷7=96, 6"; new="408÷7=58, 2"},
    @{old="801÷2=400, 1"; new="382÷4=95, 2"},
    @{old="612÷6=102, 0"; new="957÷9=106, 3"},
    @{old="304÷5=60, 4"; new="508÷2=254, 0"},
    @{old="725÷3=241, 2"; new="253÷4=63, 1"},
    @{old="404÷5=80, 4"; new="577÷3=192, 1"},
    @{old="261÷2=130, 1"; new="234÷6=39, 0"},
    @{old="479÷4=119, 3"; new="462÷9=51, 3"},
    @{old="814÷3=271, 1"; new="878÷3=292, 2"},
    @{old="788÷6=131, 2"; new="154÷4=38, 2"},
    @{old="437÷7=62, 3"; new="926÷8=115, 6"},
    @{old="962÷8=120, 2"; new="883÷6=147, 1"},
    @{old="629÷2=314, 1"; new="514÷2=257, 0"},
    @{old="595÷6=99, 1"; new="523÷7=74, 5"},
    @{old="461÷6=76, 5"; new="453÷8=56, 5"},
    @{old="746÷5=149, 1"; new="454÷7=64, 6"},
    @{old="972÷4=243, 0"; new="297÷2=148, 1"},
    @{old="995÷3=331, 2"; new="815÷6=135, 5"},
    @{old="140÷7=20, 0"; new="755÷3=251, 2"},
    @{old="776÷3=258, 2"; new="956÷2=478, 0"},
    @{old="361÷8=45, 1"; new="950÷2=475, 0"},
    @{old="878÷4=219, 2"; new="770÷9=85, 5"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
